$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure these cells keep their original "text" cell type (not auto-converted
# to numbers) by applying a text number format before assigning values that
# look numeric, matching the workbook source (inline strings).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.498.44'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -1.02%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.913.31'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -1.49%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '239.58'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -1.34%  '
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.16%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4792'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -1.88%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2843'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -3.57%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06703'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -2.83%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '18.88'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -2.92%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '102.18'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -3.64%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07713'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.25%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.917.30'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -1.26%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.207'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -2.84%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6706'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -4.08%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '269.71'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -1.28%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.509.01'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.002'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.10%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007477'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -3.23%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.68'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -3.32%  '
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -1.69%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.002'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.10%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.300'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -4.11%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.390'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -3.44%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '167.24'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.01%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '19.24'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -2.08%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.062'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -4.87%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.03%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.1002'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -3.97%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.624'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.91%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.517'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -2.40%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.215'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -3.72%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04724'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -2.85%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7260'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -3.87%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.107'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -4.41%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.723'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.40%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01916'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -4.11%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.611'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -1.79%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.296'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -3.98%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '74.79'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -4.02%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.970'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -6.20%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8598'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -4.96%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '105.22'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -2.68%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.4264'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -3.26%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.002'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.33%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.403'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -4.64%  '
$ws.Range("B47").Value = 'Algorand'
$ws.Range("C47").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.1200'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -3.81%  '
$ws.Range("B48").Value = 'Maker'
$ws.Range("C48").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '918.74'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -8.57%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '34.80'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -3.44%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.801'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -5.31%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05764'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.51%  '
